$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 196
$ws.Range("J2").Value = 710
$ws.Range("L2").Value = 198
$ws.Range("M2").Value = 8
$ws.Range("N2").Value = 121
$ws.Range("O2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 9
$ws.Range("S2").Value = 71
$ws.Range("T2").Value = 107
$ws.Range("U2").Value = 5
$ws.Range("V2").Value = 1047
$ws.Range("X2").Value = 1054
$ws.Range("Z2").Value = 13
$ws.Range("AA2").Value = 5
